# hw - ledboard, compensation
#
# FB_devider.xlsx / Sheet1: rework the feedback-divider sheet so B11 (R2)
# is solved from a compensation formula that accounts for the diode
# threshold drop at C10, instead of the old straight-line estimate. Also
# refreshes the min/max V_out seed values (B2/B3) and the diode-drop
# measurement (C10), and reflows the dependent cells (B9, C9, D9, D10,
# D11, C13) that sit downstream of those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3/B4 originally carry a cell style (fillId=2, "highlight") that is a
# byte-for-byte duplicate of the style already used by B1/B2. Re-apply
# that shared look by format-copying from B1 so the workbook doesn't keep
# growing the style table with redundant entries; C10 picks up the same
# highlight now that it is promoted to an input cell.
$ws.Range("B1").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C10").PasteSpecial(-4122)

# --- Updated inputs ---
$ws.Range("B2").Value = 0.40000000000000002
$ws.Range("B3").Value = 8.5589999999999993
$ws.Range("C10").Value = 0.59999999999999998

# --- Formula changes ---
# D10 and D11 are now one shared formula (C10/B10, shifted to C11/B11 for
# row 11) instead of two independently authored formulas.
$ws.Range("D10:D11").Formula = "=C10/B10"

# B11 (R2) is now solved with compensation for the diode drop C10.
$ws.Range("B11").Formula = "=(G1-B2-C10)/(D10-(B3-C10)/(B9))"

# Recalculate so every dependent cell (C9, D9, C11, D11, C13, ...) caches
# a fresh value.
$wb.Application.Calculate()
